# "Volviendo al inicio desde obtener excel"
#
# Renames the two placeholder contestants ("a"/"b") to their real names
# (Andrew / Daniela), adds a POSICION/RONDA label column, clears the
# per-row contestant tag on the "Puntuaciones generales" answer log and
# appends a second contestant's answers to the last 4 questions there,
# and drops the now-duplicated empty "Play" row on "Puntuaciones PASO".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Tablero 1" (ranking board)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B1").Value = "POSICION"
$ws1.Range("A2").Value = "Andrew"
$ws1.Range("A3").Value = "Daniela"

# ---------------------------------------------------------------------
# Sheet 2: "Puntuaciones generales" (full answer log)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B1").Value = "RONDA"

# Clear the per-row contestant tag in column A for every existing data row.
$ws2.Range("A2:A17").Value = ""

# Append a second contestant's answers to the last 4 "Segunda Ronda
# Grupal" / GENERALES questions (rows 18-21).
$ws2.Range("A18").Value = ""
$ws2.Range("B18").Value = "Segunda Ronda Grupal"
$ws2.Range("C18").Value = "GENERALES"
$ws2.Range("D18").Value = "Estudia las fuerzas tectónicas y procesos del interior de la Tierra. Puede estudiar las transformaciones de la estructura interna de la Tierra y/o los procesos exógenos de la superficie terrestre:"
$ws2.Range("E18").Value = "Sismotectonica"
$ws2.Range("F18").Value = -1

$ws2.Range("A19").Value = ""
$ws2.Range("B19").Value = "Segunda Ronda Grupal"
$ws2.Range("C19").Value = "GENERALES"
$ws2.Range("D19").Value = "Tipo de hidrocarburo predominante en la Cuenca de Veracruz:"
$ws2.Range("E19").Value = ""
$ws2.Range("F19").Value = -1

$ws2.Range("A20").Value = ""
$ws2.Range("B20").Value = "Segunda Ronda Grupal"
$ws2.Range("C20").Value = "GENERALES"
$ws2.Range("D20").Value = "Procedimiento de simulación y análisis estadístico que incorpora parámetros de incertidumbre para evaluar volumetría y riesgo y jerarquizar los prospectos exploratorios:"
$ws2.Range("E20").Value = ""
$ws2.Range("F20").Value = -2

$ws2.Range("A21").Value = ""
$ws2.Range("B21").Value = "Segunda Ronda Grupal"
$ws2.Range("C21").Value = "GENERALES"
$ws2.Range("D21").Value = "El ""Play"" es un término que define:"
$ws2.Range("E21").Value = "Un objetivo exploratorio"
$ws2.Range("F21").Value = 3

# ---------------------------------------------------------------------
# Sheet 3: "Puntuaciones PASO" (per-contestant final scores)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Drop the duplicated, answer-less "Play" row (old row 5); the row below
# it (the one that actually carries the answer) shifts up to take its
# place.
$ws3.Rows.Item(5).Delete()

$ws3.Range("B1").Value = "RONDA"
$ws3.Range("A2").Value = "Andrew"
$ws3.Range("A3").Value = "Daniela"
$ws3.Range("A4").Value = "Daniela"
$ws3.Range("A5").Value = "Daniela"
